# feat: add 2022-Q3 data
#
# The workbook tracks a fund's quarterly holdings: one worksheet per
# quarter (fund code / name / size / position / market value / rank) plus
# a "总计" (totals) roll-up sheet. A new quarter (2022-Q3) is being added:
#   - a new worksheet named "2022-Q3" is inserted right after "总计" (i.e.
#     right before the existing "2022-Q2" tab), carrying the new quarter's
#     figures for the same two funds shown on the other quarterly tabs;
#   - the "总计" roll-up sheet gets a new top data row for 2022-Q3; the
#     rows that were already there (2022-Q2 / 2022-Q1 / 2021-Q4) shift down
#     one row, unchanged otherwise.
# The other quarterly tabs (2022-Q2 / 2022-Q1 / 2021-Q4) keep their own
# names and figures as they were -- inserting a new tab only changes their
# tab position, never their content.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# --- 1. Insert the new "2022-Q3" sheet right after "总计" -----------------
# Duplicate the existing "2022-Q2" tab (same headers / fund codes / fund
# names / number formats) and drop the copy in right before it, then
# rename the copy and overwrite just the figures that differ for Q3.
$q2Sheet.Copy($q2Sheet, $null)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Row 2: 招商港股通核心精选股票A (011651) -- fund code/name/rank unchanged,
# only size / position / market-value / weight change for the new quarter.
$q3Sheet.Cells.Item(2, 4).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 4).Value = "1.66"
$q3Sheet.Cells.Item(2, 5).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 5).Value = "84.28"
$q3Sheet.Cells.Item(2, 6).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 6).Value = "3.74"
$q3Sheet.Cells.Item(2, 7).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 7).Value = "0.0621"

# Row 3: 招商港股通核心精选股票C (011652)
$q3Sheet.Cells.Item(3, 4).NumberFormat = "@"
$q3Sheet.Cells.Item(3, 4).Value = "0.88"
$q3Sheet.Cells.Item(3, 5).NumberFormat = "@"
$q3Sheet.Cells.Item(3, 5).Value = "84.28"
$q3Sheet.Cells.Item(3, 6).NumberFormat = "@"
$q3Sheet.Cells.Item(3, 6).Value = "3.74"
$q3Sheet.Cells.Item(3, 7).NumberFormat = "@"
$q3Sheet.Cells.Item(3, 7).Value = "0.0329"

# --- 2. Update the "总计" roll-up sheet: insert a new top row for 2022-Q3 -
# Shift the existing three data rows (rows 2-4) down to rows 3-5, then
# write the new 2022-Q3 row into row 2. Only the running index in column A
# needs correcting on the shifted rows -- the date/count/value columns
# move down untouched.
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.1

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
